$d = $word.ActiveDocument

# 1. Trim the first info line down to just the name; drop gender/ethnicity/hometown
#    (gender gets re-added below on the address line).
$d.Content.Find.Execute(
    "姓名：王美姣  性别：女   民族：汉族   籍贯：湖南",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "姓名：王美姣", 2) | Out-Null

# 2. Remove the blank line, the phone-number line, and the following blank line
#    entirely (paragraphs 3-5 once step 1 has run).
$startRange = $d.Paragraphs.Item(3).Range.Start
$endRange   = $d.Paragraphs.Item(5).Range.End
$d.Range($startRange, $endRange).Delete()

# 3. Replace the address line's text with the gender info (keeps that
#    paragraph's formatting and the _GoBack bookmark intact).
$d.Content.Find.Execute(
    "地址：湖南省益阳市",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "性别：女", 2) | Out-Null
